# "added another operation in TC 4 and added the TC 5"
# TC4 lives on the second sheet; it is renamed from "Sheet2" to "webshop" and
# gets a login-style username/password pair whose values are auto-hyperlinked
# (mirrors Excel's "smart" hyperlink-on-paste behaviour), giving us the new
# shared strings + the Hyperlink cell style + the two hyperlink relationships
# for TC5.

$wb = $excel.ActiveWorkbook

# --- rename Sheet2 -> webshop (TC5 worksheet) ---------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Name = "webshop"

# --- header row (reuses the existing "UserName"/"Password" shared strings) ---
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"

# --- data row: credentials, each one a hyperlink ------------------------
$ws.Range("A2").Value = "abhinavrevu16@gmail.com"
$ws.Range("B2").Value = "Abhi@123"

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:abhinavrevu16@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Abhi@123")

# --- column widths (autofit-ish, matches the saved workbook) ------------
$ws.Columns.Item(1).ColumnWidth = 23
$ws.Columns.Item(2).ColumnWidth = 9.25

# --- make webshop the active tab / selection -----------------------------
$ws.Activate() | Out-Null
$ws.Range("B5").Select() | Out-Null
